# Add code to tackle different unit in channel / Add LFP as well
# Fill in the "Exp_collection" (column E) as RedDim_Evol for the last four
# rows, and populate the associated Expi / pref_chan / stim_size columns
# (F/G/H) that were still missing for those rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 140
$ws.Range("E140").Value = "RedDim_Evol"
$ws.Range("H140").Value = 0.75

# Row 141
$ws.Range("E141").Value = "RedDim_Evol"
$ws.Range("H141").Value = 2

# Row 142
$ws.Range("E142").Value = "RedDim_Evol"
$ws.Range("F142").Value = 5
$ws.Range("G142").Value = 5
$ws.Range("H142").Value = 3

# Row 143
$ws.Range("E143").Value = "RedDim_Evol"
$ws.Range("F143").Value = 6
$ws.Range("G143").Value = 20
$ws.Range("H143").Value = 3

# Reflect the view state left behind by the author (scrolled down and
# selected J141) when they saved the workbook.
$win = $excel.ActiveWindow
$win.ScrollRow = 115
$win.ScrollColumn = 1
$ws.Range("J141").Select()
